$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4: wrap the trailing warranty-details image link in an <em> block with an
# <a href="..."> anchor (link text "링크") instead of the bare URL.
$d4 = $ws.Range("D4")
$oldTail = '<br><em>https://www.volvocars.com/images/v/-/media/market-assets/korea/applications/localpages/test/warranty-program/warranty_details_2.png</em>'
$newTail = '<br><em><a href="https://www.volvocars.com/images/v/-/media/market-assets/korea/applications/localpages/test/warranty-program/warranty_details_2.png">링크</a></em>'
$d4Text = $d4.Value()
$d4.Value = $d4Text.Replace($oldTail, $newTail)

# F4: the cell used to hold the bare URL string (now inlined into D4 above)
# and a hyperlink pointing at it. Remove the hyperlink and empty the cell,
# keeping its (hyperlink-styled) formatting.
$f4 = $ws.Range("F4")
$f4.Hyperlinks.Delete()
$f4.ClearContents()

# Selection moves from F4 to D4.
$ws.Range("D4").Select()
